$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2204.2292171065692
$ws.Range("B1").Value = 1384.9784196034625
$ws.Range("C1").Value = 1432.3647031410342
$ws.Range("A2").Value = 2227.1217362298958
$ws.Range("B2").Value = 1437.8378068311729
$ws.Range("C2").Value = 1597.6623857966588
$ws.Range("A3").Value = 2343.9137069389285
$ws.Range("B3").Value = 1562.5519803917387
$ws.Range("C3").Value = 1546.0202461918882
$ws.Range("A4").Value = 2320.9625107300044
$ws.Range("B4").Value = 1777.2019743706417
$ws.Range("C4").Value = 1719.7151079426428
$ws.Range("A5").Value = 2422.6539933914073
$ws.Range("B5").Value = 1667.3574578925873
$ws.Range("C5").Value = 1629.5498009388937
$ws.Range("A6").Value = 2360.708401255838
$ws.Range("B6").Value = 1774.6928232753951
$ws.Range("C6").Value = 1784.4776913416292
$ws.Range("A7").Value = 1992.9848978898274
$ws.Range("B7").Value = 1566.2710828612082
$ws.Range("C7").Value = 1484.2211590598711
$ws.Range("A8").Value = 2135.7478628526442
$ws.Range("B8").Value = 1655.7390018863714
$ws.Range("C8").Value = 1507.3584775010195
$ws.Range("A9").Value = 2471.0177526676862
$ws.Range("B9").Value = 1800.6029858790894
$ws.Range("C9").Value = 1551.7456437379697
$ws.Range("A10").Value = 2111.6060965539632
$ws.Range("B10").Value = 1384.7449821684261
$ws.Range("C10").Value = 1362.2270035448807
$ws.Range("A11").Value = 1970.0038839174615
$ws.Range("B11").Value = 1416.2529893939627
$ws.Range("C11").Value = 1298.4013325128637
$ws.Range("A12").Value = 2787.9537061107872
$ws.Range("B12").Value = 2270.4411671699713
$ws.Range("C12").Value = 2036.7399429835384
$ws.Range("A13").Value = 2315.185341097801
$ws.Range("B13").Value = 1763.908795922524
$ws.Range("C13").Value = 1803.6381877513845
$ws.Range("A14").Value = 2593.0481325924952
$ws.Range("B14").Value = 1919.8514062571228
$ws.Range("C14").Value = 1771.5022962243368
$ws.Range("A15").Value = 2508.3211273020761
$ws.Range("B15").Value = 2032.9295184874288
$ws.Range("C15").Value = 1966.7636030502799
$ws.Range("A16").Value = 2205.5211166496383
$ws.Range("B16").Value = 1536.4424786824784
$ws.Range("C16").Value = 1270.3318538283795
$ws.Range("A17").Value = 2225.5448111673422
$ws.Range("B17").Value = 1678.4907154266741
$ws.Range("C17").Value = 1588.4518536184632
$ws.Range("A18").Value = 2487.9072684430193
$ws.Range("B18").Value = 2061.2592403222206
$ws.Range("C18").Value = 1916.9785060868273
$ws.Range("A19").Value = 1810.972624286753
$ws.Range("B19").Value = 1940.946478028055
$ws.Range("C19").Value = 1955.4582997565253
$ws.Range("A20").Value = 2351.2568365264019
$ws.Range("B20").Value = 1850.0736672655316
$ws.Range("C20").Value = 1653.3462737461939
$ws.Range("A21").Value = 2583.110078207194
$ws.Range("B21").Value = 1907.0371173723911
$ws.Range("C21").Value = 1813.052769841026
$ws.Range("A22").Value = 2448.3691807930336
$ws.Range("B22").Value = 1890.2457209890522
$ws.Range("C22").Value = 1644.281071491549
